# Update crypto price/volume data per GitHub Actions scrape refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.265.26"
$ws.Range("E2").Value = "  -1.62%  "
$ws.Range("D3").Value = "2.345.17"
$ws.Range("E3").Value = "  +3.72%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("E5").Value = "  +2.02%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "230.96"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.19%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "65.28"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.62%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.458"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.46%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0946"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.04%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "56.92"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.39%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "26.62"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.43%  "
$ws.Range("D13").Value = "2.694.22"
$ws.Range("E13").Value = "  +3.59%  "
$ws.Range("E14").Value = "  -1.36%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.30"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.73%  "
$ws.Range("E16").Value = "  +2.62%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.838"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.28%  "
$ws.Range("D18").Value = "2.341.53"
$ws.Range("E18").Value = "  +3.26%  "
$ws.Range("D19").Value = "43.271.63"
$ws.Range("E19").Value = "  -1.47%  "
$ws.Range("E20").Value = "  -2.82%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "73.55"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.08%  "
$ws.Range("E22").Value = "  +1.42%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "248.13"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.33%  "
$ws.Range("E24").Value = "  +20.36%  "
$ws.Range("E25").Value = "  +0.10%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.43"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.55%  "
$ws.Range("E27").Value = "  -1.59%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "175.23"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.44%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "22.20"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +6.21%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.49"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +7.35%  "
$ws.Range("E32").Value = "  -7.47%  "
$ws.Range("E33").Value = "  +0.33%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.97"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.06%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0685"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.86%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.98"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.58%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.47"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +7.06%  "
$ws.Range("E38").Value = "  -0.68%  "
$ws.Range("E39").Value = "  -5.71%  "
$ws.Range("E40").Value = "  -2.74%  "
$ws.Range("E41").Value = "  +0.02%  "
$ws.Range("E42").Value = "  +7.89%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "17.84"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.76%  "
$ws.Range("E44").Value = "  +7.11%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "98.36"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.23%  "
$ws.Range("E46").Value = "  -0.25%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.37"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.41%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0944"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.41%  "
$ws.Range("D49").Value = "1.434.45"
$ws.Range("E49").Value = "  -0.60%  "
$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D50").Value = "2.567.45"
$ws.Range("E50").Value = "  +3.78%  "
$ws.Range("B51").Value = "TerraClassic"
$ws.Range("C51").Value = "https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.000203"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -8.54%  "
